# all_args_here small print bug fixed in on.exit()
#
# The "all_args_here" row was missing its "x" mark in column D (the
# print was silently dropped because of a bug in on.exit()). This adds
# the missing "x" to D2, matching the formatting already used by the
# neighbouring B2/C2 cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of C2 (centered, wrapped, bold header style)
# onto D2, then fill in the same "x" mark used throughout column B/C/D.
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D2").Value2 = $ws.Range("C2").Value2

# Leave the active cell/selection on D2, as it was the last edited cell.
$ws.Range("D2").Select()
